$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9159573316574097
$ws.Range("B1").Value = 1.420223236083984
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 1.402090072631836
$ws.Range("E1").Value = 1.392821669578552
